# Remove unnecessary source_id field from the survey sheet.
# The source_id field lives on row 5 of the "survey" sheet
# (type=text, name=source_id, label=Source ID, appearance=hidden).
# Deleting the whole row shifts everything below it up by one,
# which matches the target edit exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

$ws.Rows.Item(5).Delete()
